$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data one column to the left, overwriting the "split" column (C)
# with what used to be in D:M (TV, Digital, ... Radio, Other). This is a
# plain copy/paste (not cut), so the source column M is left with its
# original values untouched.
$ws.Range("D1:M54").Copy()
$ws.Range("C1").PasteSpecial()

# Renumber the "week" column (B) for the 2016 rows (originally rows 29-54)
# to a sequential 0-based count instead of continuing from the 2015 weeks.
for ($i = 0; $i -le 25; $i++) {
    $row = 29 + $i
    $ws.Cells.Item($row, 2).Value = $i
}

# Update the sheet view selection to D9
$ws.Range("D9").Select()
